$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("K6").Value = -0.162895509468913
$ws.Range("E12").Value = -0.162895509468913

$ws.Range("U8").Value = 0.175805941336879
$ws.Range("G22").Value = 0.175805941336879

$ws.Range("N11").Value = 0.18060041601226
$ws.Range("J15").Value = 0.18060041601226

$ws.Range("N13").Value = -0.174109457732615
$ws.Range("L15").Value = -0.174109457732615
